$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.071.57"
$ws.Range("E2").Value = "  -0.79%  "

$ws.Range("D3").Value = "3.403.64"
$ws.Range("E3").Value = "  -1.16%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "'572.69"
$ws.Range("E5").Value = "  -0.09%  "

$ws.Range("D6").Value = "'162.51"
$ws.Range("E6").Value = "  +2.29%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").Value = "3.401.46"
$ws.Range("E8").Value = "  -1.14%  "

$ws.Range("D9").Value = "'0.550"
$ws.Range("E9").Value = "  -5.23%  "

$ws.Range("E10").Value = "  +1.26%  "

$ws.Range("E11").Value = "  -1.98%  "

$ws.Range("E12").Value = "  -4.37%  "

$ws.Range("D13").Value = "3.992.04"
$ws.Range("E13").Value = "  -1.10%  "

$ws.Range("E14").Value = "  +0.49%  "

$ws.Range("D15").Value = "'26.85"
$ws.Range("E15").Value = "  -2.25%  "

$ws.Range("E16").Value = "  -0.72%  "

$ws.Range("D17").Value = "64.103.51"
$ws.Range("E17").Value = "  -0.77%  "

$ws.Range("D18").Value = "3.440.98"
$ws.Range("E18").Value = "  -0.20%  "

$ws.Range("E19").Value = "  -0.98%  "

$ws.Range("D20").Value = "'13.41"
$ws.Range("E20").Value = "  -2.49%  "

$ws.Range("D21").Value = "'374.18"
$ws.Range("E21").Value = "  -1.55%  "

$ws.Range("D22").Value = "'7.77"
$ws.Range("E22").Value = "  -2.36%  "

$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  +0.19%  "

$ws.Range("E24").Value = "  -2.63%  "

$ws.Range("E25").Value = "  -3.56%  "

$ws.Range("E26").Value = "  -4.65%  "

$ws.Range("D27").Value = "'9.48"
$ws.Range("E27").Value = "  -4.50%  "

$ws.Range("E28").Value = "  -0.68%  "

$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -1.84%  "

$ws.Range("E30").Value = "  -0.55%  "

$ws.Range("E31").Value = "  -3.20%  "

$ws.Range("E32").Value = "  -0.22%  "

$ws.Range("D33").Value = "'1.00"
$ws.Range("E33").Value = "  +0.07%  "

$ws.Range("E34").Value = "  -1.94%  "

$ws.Range("D35").Value = "'7.01"
$ws.Range("E35").Value = "  +0.11%  "

$ws.Range("D36").Value = "'1.47"
$ws.Range("E36").Value = "  -6.37%  "

$ws.Range("D37").Value = "'159.37"
$ws.Range("E37").Value = "  -1.11%  "

$ws.Range("E38").Value = "  +7.33%  "

$ws.Range("E39").Value = "  -2.79%  "

$ws.Range("D40").Value = "'0.0724"
$ws.Range("E40").Value = "  -2.85%  "

$ws.Range("D41").Value = "'25.74"
$ws.Range("E41").Value = "  -2.02%  "

$ws.Range("D42").Value = "'42.63"
$ws.Range("E42").Value = "  -0.67%  "

$ws.Range("D43").Value = "2.722.50"
$ws.Range("E43").Value = "  -5.49%  "

$ws.Range("E44").Value = "  -1.09%  "

$ws.Range("D45").Value = "'25.73"
$ws.Range("E45").Value = "  -0.28%  "

$ws.Range("D46").Value = "'4.34"
$ws.Range("E46").Value = "  -3.98%  "

$ws.Range("E47").Value = "  -1.98%  "

$ws.Range("D48").Value = "'2.39"
$ws.Range("E48").Value = "  -1.78%  "

$ws.Range("D49").Value = "'328.55"
$ws.Range("E49").Value = "  +1.86%  "

$ws.Range("E50").Value = "  -2.66%  "

$ws.Range("E51").Value = "  -1.99%  "
